$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 becomes the "end of group" row: apply the bordered style (like rows 3/6/10) ---
$ws.Range("A3:B3").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C3:E3").Copy()
$ws.Range("C12:E12").PasteSpecial(-4122) # xlPasteFormats

# --- New row 13: same (unbordered) style as row 11 ---
$ws.Range("B11:E11").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122) # xlPasteFormats

# --- New row 14: same (unbordered) style as row 9 ---
$ws.Range("B9:E9").Copy()
$ws.Range("B14:E14").PasteSpecial(-4122) # xlPasteFormats

# --- Values for row 13 ---
$ws.Range("B13").Value = 143
$ws.Range("C13").Value = " Congratulations on graduating!\nYou two are extraordinary!"

# --- Values for row 14 ---
$ws.Range("B14").Value = 146
$ws.Range("C14").Value = " Do you know what they say\nabout that guild?[K] That they make it terribly\nhard for apprentices to graduate!"

# --- Remaining translated / converted strings (order matches the authoring tool's append order) ---
$ws.Range("D13").Value = " Поздравляю с выпуском! Вы двое\nневероятны!"
$ws.Range("E13").Value = " Ðïèäñàâìÿý ò âúðôòëïí! Âú äâïå\nîåâåñïÿóîú!"
$ws.Range("D14").Value = " Знаете, что говорят о гильдии?[K]\nТо, что ученикам ужасно тяжело выпуститься\nиз неё!"
$ws.Range("E14").Value = " Èîàåóå, œóï ãïâïñÿó ï ãéìûäéé?[K]\nÓï, œóï ôœåîéëàí ôçàòîï óÿçåìï âúðôòóéóûòÿ\néè îåæ!"

# --- Row heights for the two new rows ---
$ws.Rows.Item(13).RowHeight = 21.6
$ws.Rows.Item(14).RowHeight = 42

# --- Update selection to match the new active cell ---
$ws.Range("E14").Select()
